$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2023-11-30 Thursday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-12-01 Friday", 2)

# Update the division practice table. Addressing cells directly by
# (row, column) avoids any ambiguity from duplicate cell text appearing
# elsewhere in the table.
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    $r.End = $r.End - 1
    $r.Text = $newText
}

Set-CellText $t 1 1 "64÷9=7, 1"
Set-CellText $t 1 2 "98÷6=16, 2"
Set-CellText $t 1 3 "11÷6=1, 5"
Set-CellText $t 1 4 "39÷7=5, 4"
Set-CellText $t 1 5 "24÷2=12, 0"

Set-CellText $t 5 1 "52÷8=6, 4"
Set-CellText $t 5 2 "50÷9=5, 5"
Set-CellText $t 5 3 "64÷3=21, 1"
Set-CellText $t 5 4 "29÷3=9, 2"
Set-CellText $t 5 5 "88÷8=11, 0"

Set-CellText $t 9 1 "48÷6=8, 0"
Set-CellText $t 9 2 "70÷8=8, 6"
Set-CellText $t 9 3 "71÷2=35, 1"
Set-CellText $t 9 4 "53÷8=6, 5"
Set-CellText $t 9 5 "90÷6=15, 0"

Set-CellText $t 13 1 "87÷6=14, 3"
Set-CellText $t 13 2 "49÷5=9, 4"
Set-CellText $t 13 3 "67÷4=16, 3"
Set-CellText $t 13 4 "74÷8=9, 2"
Set-CellText $t 13 5 "56÷6=9, 2"

Set-CellText $t 17 1 "18÷3=6, 0"
Set-CellText $t 17 2 "86÷9=9, 5"
Set-CellText $t 17 3 "99÷8=12, 3"
Set-CellText $t 17 4 "58÷9=6, 4"
Set-CellText $t 17 5 "92÷9=10, 2"
